$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (13 and 14) to hold the responsible-teacher entries;
# everything from the old row 13 downward shifts down by two rows.
$ws.Rows("13:14").Insert()

# The inserted rows pick up column-default formatting which does not match
# the wrap/top-aligned "data" style used throughout columns B/C, so copy the
# correct format over from a row that already has it before filling values.
$ws.Range("B15:C15").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B15:C15").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Objetivos: (row 10) now holds the real Portuguese objectives text ---
$objPt = @'
Apresentação e aplicação dos fundamentos teóricos das operações unitárias envolvendo transferência de calor e massa. Os tópicos abordados constituem aplicação prática dos conhecimentos desenvolvidos ao longo da disciplina fenômenos de transporte II e são de grande importância para estudos posteriores de processos químicos industriais.
'@
$ws.Range("B10").Value = $objPt
$ws.Range("C10").Value = $objPt

# --- Docentes responsáveis (rows 13-14, newly inserted) ---
$ws.Range("B13").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("C13").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("B14").Value = "5817372 - Simone de Fátima Medeiros Sampaio"
$ws.Range("C14").Value = "5817372 - Simone de Fátima Medeiros Sampaio"

# --- Programa resumido: (row 15, was row 13) ---
$progResumidoPt = @'
1)Trocadores de Calor;
2)Trocadores de Calor Tubulares;
3)Trocadores de calor de Placas;
4)Sistemas de Troca de Calor com mudança de Fase;
5)Psicrometria;
6)Umidificação e Desumidificação;
'@
$ws.Range("B15").Value = $progResumidoPt
$ws.Range("C15").Value = $progResumidoPt

# --- Programa: (row 17, was row 15) full detailed syllabus ---
$programaPt = @'
1)Trocadores de Calor: conceitos gerais e tipos de trocadores de calor;
2)Trocadores de calor tubulares: cálculos em um trocador de calor bitubular;  Método DTML; Trocadores de calor de casco e tubos; Correlações para determinação dos coeficientes de transferência de calor em trocadores de calor de casco e tubos; Estimativa dos coeficientes de película; Método ε-NUT; Queda de pressão nos trocadores de casco e tubos;
3)Trocadores de calor de placas: cálculos e comparação com trocadores tubulares;
4)Sistemas de troca de calor com mudança de fase: evaporadores, condensadores , refervedores e caldeiras; Cristalização;
5)Psicrometria: conceitos envolvidos e uso da carta psicrométrica;
6)Operações de umidificação de desumidificação; Torres de resfriamento e Secagem.
'@
$ws.Range("B17").Value = $programaPt
$ws.Range("C17").Value = $programaPt

# --- Método: (row 20, was row 18) ---
$ws.Range("B20").Value = "Aplicação de 2 provas, P1 e P2."
$ws.Range("C20").Value = "Aplicação de 2 provas, P1 e P2."

# --- Critério: (row 21, was row 19) ---
$criterioB = @'
A média do período será MP = (P1+P2)/2. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham frequência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou frequência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham frequência mínima de 70% serão submetidos ao período de recuperação (regimental).
'@
$ws.Range("B21").Value = $criterioB
$ws.Range("C21").Value = $criterioB

# --- Norma de recuperação: (row 22, was row 20) ---
$normaB = @'
A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação.
'@
$ws.Range("B22").Value = $normaB
$ws.Range("C22").Value = $normaB

# --- Bibliografia: (row 23, was row 21) ---
$biblioB = @'
1)COULSON, J. M.; RICHARDSON; J.F. Chemical Engineering. v. 2: Particle Technology e Separation Processes. 5ed. Amsterdan: Butterworth Heinemann, 1229p. 2005;
2)COULSON & Richardson's Chemical Engineering: chemical engineering design by R.K. Sinnott. 6ed. Amsterdam: Elsevier Butterworth Heinemann, 895p. 2004;
3)COUPER, J. R.; PENNEY, W. R.; FAIR, J. R.; W.; Stanley. M. Chemical Process Equipment: Selection and Design. 2ed. Amsterdam: Elsevier, 814p. 2005;
4)FOUST, A. S.; WENZEL, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSEN, L. B. 2ed. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 670p. 2008;
5)GEANKOPLIS, C. J. Transport Processes and Separation Process Principles. 4ed. New York: Prentice Hall, 1026p. 2010;
6)HIMMELBLAU, D. M.; RIGGS, J. B. Engenharia Química: princípios e cálculos. 7ed. Rio de Janeiro: LTC, 846p. 2006;
7)KERN, D. Q. Processos de transmissão de calor. Rio de Janeiro: Guanabara Dois, 671p. 1982;
8)MCCABE, W. L.; SMITH, J. C.; HARRIOT, P. Unit operations of chemical engineering. 7ed. Boston: McGraw-Hill, 1140 p. 2005;
9)PERRY's chemical engineers handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry New York: McGraw-Hill, 2008.
'@
$ws.Range("B23").Value = $biblioB
$ws.Range("C23").Value = $biblioB
